# Updates the "cryptos" price/volume table (columns D and E, rows 2-51)
# with freshly scraped values, per the Sep 1 2023 GitHub Actions run.
#
# Column D ("Price") cells are stored as TEXT in the workbook
# (e.g. "26.116.84", "1.008") even though some of the new values are
# syntactically valid numbers (e.g. "1.008", "215.63"). Writing a
# numeric-looking string straight into Range.Value would let Excel
# auto-convert it to a real number, changing the cell's stored type.
# To keep those specific cells as text (matching the original
# authoring), each update row below carries an 'AsText' flag
# (precomputed from the target value) that forces the cell's
# NumberFormat to Text ("@") immediately before the value is written.
# Values that Excel would never mistake for a number (two '.'s,
# subscript digits, the '%' volume strings, ...) leave AsText = $false
# and are written as-is, exactly like the rest of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '26.116.84'; AsText = $false },
    @{ Cell = 'E2'; Value = '  -4.50%  '; AsText = $false },
    @{ Cell = 'D3'; Value = '1.650.69'; AsText = $false },
    @{ Cell = 'E3'; Value = '  -3.62%  '; AsText = $false },
    @{ Cell = 'D4'; Value = '1.008'; AsText = $true },
    @{ Cell = 'E4'; Value = '  -0.06%  '; AsText = $false },
    @{ Cell = 'D5'; Value = '215.63'; AsText = $true },
    @{ Cell = 'E5'; Value = '  -3.92%  '; AsText = $false },
    @{ Cell = 'D6'; Value = '0.5112'; AsText = $true },
    @{ Cell = 'E6'; Value = '  -2.86%  '; AsText = $false },
    @{ Cell = 'D7'; Value = '1.008'; AsText = $true },
    @{ Cell = 'E7'; Value = '  +0.01%  '; AsText = $false },
    @{ Cell = 'D8'; Value = '0.2589'; AsText = $true },
    @{ Cell = 'E8'; Value = '  -1.86%  '; AsText = $false },
    @{ Cell = 'D9'; Value = '0.06436'; AsText = $true },
    @{ Cell = 'E9'; Value = '  -2.90%  '; AsText = $false },
    @{ Cell = 'D10'; Value = '19.72'; AsText = $true },
    @{ Cell = 'E10'; Value = '  -4.66%  '; AsText = $false },
    @{ Cell = 'D11'; Value = '0.07729'; AsText = $true },
    @{ Cell = 'E11'; Value = '  -0.22%  '; AsText = $false },
    @{ Cell = 'D12'; Value = '1.671.27'; AsText = $false },
    @{ Cell = 'E12'; Value = '  -2.42%  '; AsText = $false },
    @{ Cell = 'D13'; Value = '4.275'; AsText = $true },
    @{ Cell = 'E13'; Value = '  -3.98%  '; AsText = $false },
    @{ Cell = 'D14'; Value = '1.878.90'; AsText = $false },
    @{ Cell = 'D15'; Value = '0.5484'; AsText = $true },
    @{ Cell = 'E15'; Value = '  -4.86%  '; AsText = $false },
    @{ Cell = 'D16'; Value = '0.0₅8002'; AsText = $false },
    @{ Cell = 'E16'; Value = '  -1.98%  '; AsText = $false },
    @{ Cell = 'D17'; Value = '63.84'; AsText = $true },
    @{ Cell = 'E17'; Value = '  -5.59%  '; AsText = $false },
    @{ Cell = 'D18'; Value = '26.122.66'; AsText = $false },
    @{ Cell = 'E18'; Value = '  -4.57%  '; AsText = $false },
    @{ Cell = 'D19'; Value = '1.008'; AsText = $true },
    @{ Cell = 'E19'; Value = '  -0.02%  '; AsText = $false },
    @{ Cell = 'D20'; Value = '207.53'; AsText = $true },
    @{ Cell = 'E20'; Value = '  -5.16%  '; AsText = $false },
    @{ Cell = 'D21'; Value = '4.387'; AsText = $true },
    @{ Cell = 'E21'; Value = '  -5.36%  '; AsText = $false },
    @{ Cell = 'D22'; Value = '10.04'; AsText = $true },
    @{ Cell = 'E22'; Value = '  -3.53%  '; AsText = $false },
    @{ Cell = 'D23'; Value = '6.017'; AsText = $true },
    @{ Cell = 'E23'; Value = '  -0.02%  '; AsText = $false },
    @{ Cell = 'D25'; Value = '1.874'; AsText = $true },
    @{ Cell = 'E25'; Value = '  +8.52%  '; AsText = $false },
    @{ Cell = 'D26'; Value = '143.08'; AsText = $true },
    @{ Cell = 'E26'; Value = '  -1.59%  '; AsText = $false },
    @{ Cell = 'D27'; Value = '0.1169'; AsText = $true },
    @{ Cell = 'E27'; Value = '  -2.55%  '; AsText = $false },
    @{ Cell = 'D28'; Value = '6.919'; AsText = $true },
    @{ Cell = 'D29'; Value = '15.81'; AsText = $true },
    @{ Cell = 'E29'; Value = '  -1.98%  '; AsText = $false },
    @{ Cell = 'D30'; Value = '0.05082'; AsText = $true },
    @{ Cell = 'E30'; Value = '  -4.21%  '; AsText = $false },
    @{ Cell = 'D31'; Value = '1.243'; AsText = $true },
    @{ Cell = 'E31'; Value = '  -3.98%  '; AsText = $false },
    @{ Cell = 'D32'; Value = '3.348'; AsText = $true },
    @{ Cell = 'E32'; Value = '  -3.52%  '; AsText = $false },
    @{ Cell = 'D33'; Value = '3.233'; AsText = $true },
    @{ Cell = 'E33'; Value = '  -3.42%  '; AsText = $false },
    @{ Cell = 'D34'; Value = '1.550'; AsText = $true },
    @{ Cell = 'E34'; Value = '  -5.33%  '; AsText = $false },
    @{ Cell = 'D35'; Value = '2.349'; AsText = $true },
    @{ Cell = 'E35'; Value = '  -2.21%  '; AsText = $false },
    @{ Cell = 'D36'; Value = '0.9160'; AsText = $true },
    @{ Cell = 'E36'; Value = '  -3.42%  '; AsText = $false },
    @{ Cell = 'D37'; Value = '2.635'; AsText = $true },
    @{ Cell = 'E37'; Value = '  -7.04%  '; AsText = $false },
    @{ Cell = 'D38'; Value = '0.5696'; AsText = $true },
    @{ Cell = 'E38'; Value = '  -2.85%  '; AsText = $false },
    @{ Cell = 'D39'; Value = '1.145.06'; AsText = $false },
    @{ Cell = 'E39'; Value = '  -3.52%  '; AsText = $false },
    @{ Cell = 'D40'; Value = '0.01574'; AsText = $true },
    @{ Cell = 'E40'; Value = '  -4.52%  '; AsText = $false },
    @{ Cell = 'E41'; Value = '  +0.13%  '; AsText = $false },
    @{ Cell = 'D42'; Value = '1.008'; AsText = $true },
    @{ Cell = 'E42'; Value = '  +0.00%  '; AsText = $false },
    @{ Cell = 'D43'; Value = '5.668'; AsText = $true },
    @{ Cell = 'E43'; Value = '  -2.01%  '; AsText = $false },
    @{ Cell = 'D44'; Value = '0.8251'; AsText = $true },
    @{ Cell = 'E44'; Value = '  -1.69%  '; AsText = $false },
    @{ Cell = 'D45'; Value = '100.11'; AsText = $true },
    @{ Cell = 'E45'; Value = '  -0.93%  '; AsText = $false },
    @{ Cell = 'D46'; Value = '1.789.95'; AsText = $false },
    @{ Cell = 'E46'; Value = '  -3.64%  '; AsText = $false },
    @{ Cell = 'D47'; Value = '0.0₈112'; AsText = $false },
    @{ Cell = 'E47'; Value = '  -4.86%  '; AsText = $false },
    @{ Cell = 'D48'; Value = '0.4535'; AsText = $true },
    @{ Cell = 'E48'; Value = '  -0.58%  '; AsText = $false },
    @{ Cell = 'E49'; Value = '  +0.59%  '; AsText = $false },
    @{ Cell = 'D50'; Value = '55.22'; AsText = $true },
    @{ Cell = 'E50'; Value = '  -3.78%  '; AsText = $false },
    @{ Cell = 'D51'; Value = '7.818'; AsText = $true },
    @{ Cell = 'E51'; Value = '  -3.92%  '; AsText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)

    if ($u.AsText) {
        $cell.NumberFormat = "@"
    }

    $cell.Value = $u.Value
}
